# forests-scraped.xlsx update - 2025-10-03 12:16
# Moves the 3 rows that were on "New" into "Previously added" (appended at
# the bottom) and replaces them on "New" with 6 freshly scraped listings.

$wb     = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Archive the 3 current "New" rows onto the bottom of "Previously added"
#    (rows 147-149), carrying over their existing hyperlink targets.
# ---------------------------------------------------------------------

# Stash a clean copy of the data-row format off to the side - Hyperlinks.Add
# recolors column A with the built-in "Hyperlink" font, so every row gets
# its format re-applied from this untouched donor after the link is added.
$wsPrev.Range("A146:F146").Copy() | Out-Null
$wsPrev.Range("A500:F500").PasteSpecial($xlPasteFormats) | Out-Null

# Cadastre numbers are plain digits but must stay TEXT (they're shared
# strings in the source file, not numeric cells) - a leading "'" forces
# Excel to store them as text instead of auto-converting to a number.
$archiveRows = @(
    @{ Row=147; Link="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/udrisu-pag/bkbix.html";
       Price="35 880 €"; District="Krāslava un raj."; Area="8 ha."; Cadastre="'60960020067/68";
       Date=45932.393055555556 },
    @{ Row=148; Link="https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/salacgrivas-l-t/cdkbp.html";
       Price="50 000 €"; District="Limbaži un raj."; Area="3 ha."; Cadastre="'66720040252";
       Date=45931.83125 },
    @{ Row=149; Link="https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/liepupes-pag/eepmh.html";
       Price="39 000 €"; District="Limbaži un raj."; Area="6.50 ha."; Cadastre="'66600090044";
       Date=45931.757638888885 }
)

foreach ($r in $archiveRows) {
    $row = $r.Row

    $wsPrev.Range("A500:F500").Copy() | Out-Null
    $wsPrev.Range("A" + $row + ":F" + $row).PasteSpecial($xlPasteFormats) | Out-Null

    $wsPrev.Range("A" + $row).Value = $r.Link
    $wsPrev.Range("B" + $row).Value = $r.Price
    $wsPrev.Range("C" + $row).Value = $r.District
    $wsPrev.Range("D" + $row).Value = $r.Area
    $wsPrev.Range("E" + $row).Value = $r.Cadastre
    $wsPrev.Range("F" + $row).Value = $r.Date

    $wsPrev.Hyperlinks.Add($wsPrev.Range("A" + $row), $r.Link) | Out-Null

    # Restore the sheet's own link style (matches the other 146 rows).
    $wsPrev.Range("A500:F500").Copy() | Out-Null
    $wsPrev.Range("A" + $row + ":F" + $row).PasteSpecial($xlPasteFormats) | Out-Null
}

$wsPrev.Range("A500:F500").Clear()
$wsPrev.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Clear the old "New" rows 2-4 (values + hyperlinks) and replace with
#    the 6 newly scraped listings.
# ---------------------------------------------------------------------

# Same donor-row trick as above, captured before anything is cleared.
$wsNew.Range("A4:F4").Copy() | Out-Null
$wsNew.Range("A500:F500").PasteSpecial($xlPasteFormats) | Out-Null

$wsNew.Range("A2").Hyperlinks.Delete()
$wsNew.Range("A3").Hyperlinks.Delete()
$wsNew.Range("A4").Hyperlinks.Delete()
$wsNew.Range("A2:F4").ClearContents()

$newRows = @(
    @{ Row=2; Link="https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/auru-pag/ionih.html";
       Price="50 000 €"; District="Dobele un raj."; Area="12 ha."; Cadastre="'46460080148";
       Date=45933.61111111111 },
    @{ Row=3; Link="https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cggbgp.html";
       Price="5 000 €"; District="Kuldīga un raj."; Area="2 ha."; Cadastre="'62720050064";
       Date=45933.45347222222 },
    @{ Row=4; Link="https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cggbdp.html";
       Price="26 500 €"; District="Kuldīga un raj."; Area="5.12 ha."; Cadastre="'";
       Date=45932.93611111111 },
    @{ Row=5; Link="https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cggbdm.html";
       Price="22 000 €"; District="Kuldīga un raj."; Area="4.36 ha."; Cadastre="'";
       Date=45932.929861111115 },
    @{ Row=6; Link="https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/dunikas-pag/ngxix.html";
       Price="13 000 €"; District="Liepāja un raj."; Area="3 ha."; Cadastre="'64520070098";
       Date=45933.45416666666 },
    @{ Row=7; Link="https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/aronas-pag/ahkjo.html";
       Price="155 000 €"; District="Madona un raj."; Area="27 ha."; Cadastre="'70420020025";
       Date=45933.45972222222 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $wsNew.Range("A500:F500").Copy() | Out-Null
    $wsNew.Range("A" + $row + ":F" + $row).PasteSpecial($xlPasteFormats) | Out-Null

    $wsNew.Range("A" + $row).Value = $r.Link
    $wsNew.Range("B" + $row).Value = $r.Price
    $wsNew.Range("C" + $row).Value = $r.District
    $wsNew.Range("D" + $row).Value = $r.Area
    $wsNew.Range("E" + $row).Value = $r.Cadastre
    $wsNew.Range("F" + $row).Value = $r.Date

    $wsNew.Hyperlinks.Add($wsNew.Range("A" + $row), $r.Link) | Out-Null

    # Restore the sheet's own link style (matches row 4 / the rest of the table).
    $wsNew.Range("A500:F500").Copy() | Out-Null
    $wsNew.Range("A" + $row + ":F" + $row).PasteSpecial($xlPasteFormats) | Out-Null
}

$wsNew.Range("A500:F500").Clear()
$wsNew.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Hyperlinks.Add silently registers a built-in "Hyperlink" cell style
#    the first time it's used; the sheets here use their own pre-existing
#    styles (restored above via PasteSpecial), so drop the unused one.
# ---------------------------------------------------------------------
for ($i = $wb.Styles.Count; $i -ge 1; $i--) {
    $s = $wb.Styles.Item($i)
    if ($s.Name -eq "Hyperlink") {
        $s.Delete()
    }
}
